$wb = $excel.ActiveWorkbook

# --- Metadata sheet (sheet 1): update Date, insert "Jurisdiction" row ---
$ws = $wb.Worksheets.Item(1)

# Create row 15 first by direct-copying row 14 (content + style) so the new
# row keeps the same formatting ("s=2") as the rest of the data rows.
$ws.Range("A14:B14").Copy($ws.Range("A15:B15"))

# Read current (pre-shift) values of rows 11-13 before overwriting them.
$a11 = $ws.Range("A11").Value2
$b11 = $ws.Range("B11").Value2
$a12 = $ws.Range("A12").Value2
$b12 = $ws.Range("B12").Value2
$a13 = $ws.Range("A13").Value2
$b13 = $ws.Range("B13").Value2

# Shift rows 11-13 down into 12-14 (row 14's old content already copied to 15).
$ws.Range("A14").Value = $a13
$ws.Range("B14").Value = $b13
$ws.Range("A13").Value = $a12
$ws.Range("B13").Value = $b12
$ws.Range("A12").Value = $a11
$ws.Range("B12").Value = $b11

# New row 11: Jurisdiction property with an empty value.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update the Date value.
$ws.Range("B8").Value = "2024-09-12T14:01:50+00:00"

# --- Rename the "Include from LTI Detailed Des" sheet tab ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"
